$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# Row/column -> new value, for the "bias" section value cells (column 6)
# that were previously empty placeholders.
$updates = @{
    3  = "1.797"   # AVDD [V]
    4  = "137.2"   # IVDD [mA]
    5  = "1.798"   # DVDD [V]
    6  = "5"       # IDVDD [mA]
    7  = "2.995"   # 3V3 [V]
    8  = "3"       # I3V3 [mA]
    9  = "5.00"    # Ibias [mA]
    10 = "1.198"   # VCMSH [V]
    11 = "0.891"   # VCM [V]
    12 = "0.899"   # RVCM [V]
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 6)
    $cell.Range.Text = $updates[$row]
}
